# Update the "Timestamp" column (A) of the "Output" sheet with the
# timestamps recorded for the latest test run (24/01/2022), replacing the
# previous run's timestamps (23/01/2022) that were there before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

$ws.Range("A2").Value  = "24/01/2022 11:22:15 am"
$ws.Range("A3").Value  = "24/01/2022 11:22:22 am"
$ws.Range("A4").Value  = "24/01/2022 11:22:38 am"
$ws.Range("A5").Value  = "24/01/2022 11:22:46 am"
$ws.Range("A6").Value  = "24/01/2022 11:22:59 am"
$ws.Range("A7").Value  = "24/01/2022 11:23:15 am"
$ws.Range("A8").Value  = "24/01/2022 11:23:20 am"
$ws.Range("A9").Value  = "24/01/2022 11:23:30 am"
$ws.Range("A10").Value = "24/01/2022 11:23:35 am"
$ws.Range("A11").Value = "24/01/2022 11:24:00 am"
